$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 27347.088
$ws.Range("I11").Value = 27347.088
$ws.Range("K11").Value = 27347.088
$ws.Range("M11").Value = -27207.088
$ws.Range("H62").Value = 2237.1428
$ws.Range("I62").Value = 2375
$ws.Range("K62").Value = 2375
$ws.Range("M62").Value = -1751
$ws.Range("H65").Value = 2237.1428
$ws.Range("I65").Value = 2375
$ws.Range("K65").Value = 11875
$ws.Range("M65").Value = -8755
$ws.Range("H87").Value = 73677.78
$ws.Range("I87").Value = 14500
$ws.Range("K87").Value = 14500
$ws.Range("M87").Value = -13252
$ws.Range("H90").Value = 73677.78
$ws.Range("I90").Value = 14500
$ws.Range("K90").Value = 43500
$ws.Range("M90").Value = -37260
$ws.Range("H94").Value = 1261.5714
$ws.Range("I94").Value = 1261.5714
$ws.Range("K94").Value = 1261.5714
$ws.Range("M94").Value = -810.5714
$ws.Range("H96").Value = 1303.5454
$ws.Range("I96").Value = 1144.6154
$ws.Range("J96").Value = 1533.1111
$ws.Range("K96").Value = 3433.8462
$ws.Range("L96").Value = 4599.3333
$ws.Range("M96").Value = -2060.8462
$ws.Range("N96").Value = -7345.3333
$ws.Range("H99").Value = 589.8
$ws.Range("I99").Value = 713.625
$ws.Range("K99").Value = 2140.875
$ws.Range("M99").Value = -642.875
$ws.Range("H107").Value = 2524.75
$ws.Range("J107").Value = 3499
$ws.Range("L107").Value = 3499
$ws.Range("N107").Value = -7339
$ws.Range("H111").Value = 4424.143
$ws.Range("J111").Value = 4999.75
$ws.Range("L111").Value = 14999.25
$ws.Range("N111").Value = -21133.25
$ws.Range("H119").Value = 1530
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").Value = $null
$ws.Range("H121").Value = 1842.3636
$ws.Range("J121").Value = 1842.3636
$ws.Range("L121").Value = 5527.0908
$ws.Range("N121").Value = -9021.0908
$ws.Range("H125").Value = 250003100
$ws.Range("I125").Value = 1000000000
$ws.Range("J125").Value = 4145
$ws.Range("K125").Value = 9000000000
$ws.Range("L125").Value = 37305
$ws.Range("M125").Value = -8999997540
$ws.Range("N125").Value = -42225
$ws.Range("H127").Value = 2000
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").Value = $null
$ws.Range("H45").Value = 87288.664
$ws.Range("I45").Value = 114860.445
$ws.Range("J45").Value = 4573.3335
$ws.Range("K45").Value = 114860.445
$ws.Range("L45").Value = 4573.3335
$ws.Range("M45").Value = -114483.445
$ws.Range("N45").Value = -5327.3335
$ws.Range("H50").Value = 834.5
$ws.Range("I50").Value = 97.333336
$ws.Range("J50").Value = 1276.8
$ws.Range("K50").Value = 97.333336
$ws.Range("L50").Value = 1276.8
$ws.Range("M50").Value = 616.666664
$ws.Range("N50").Value = -2704.8
$ws.Range("H132").Value = 1673.7457
$ws.Range("I132").Value = 1154.6111
$ws.Range("K132").Value = 3463.8333
$ws.Range("M132").Value = -933.8333000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 3841
$ws.Range("I5").Value = 350
$ws.Range("J5").Value = 6168.3335
$ws.Range("K5").Value = 350
$ws.Range("L5").Value = 6168.3335
$ws.Range("M5").Value = -237
$ws.Range("N5").Value = -6394.3335
$ws.Range("H22").Value = 2027.3334
$ws.Range("J22").Value = 1949
$ws.Range("L22").Value = 1949
$ws.Range("N22").Value = -2295
$ws.Range("H33").Value = 5649.625
$ws.Range("I33").Value = 2399.6667
$ws.Range("J33").Value = 7599.6
$ws.Range("K33").Value = 2399.6667
$ws.Range("L33").Value = 7599.6
$ws.Range("M33").Value = -2063.6667
$ws.Range("N33").Value = -8271.6
$ws.Range("H86").Value = 5187.222
$ws.Range("I86").Value = 2341.8
$ws.Range("J86").Value = 8744
$ws.Range("K86").Value = 2341.8
$ws.Range("L86").Value = 8744
$ws.Range("M86").Value = -1218.8
$ws.Range("N86").Value = -10990
$ws.Range("H89").Value = 5187.222
$ws.Range("I89").Value = 2341.8
$ws.Range("J89").Value = 8744
$ws.Range("K89").Value = 11709
$ws.Range("L89").Value = 43720
$ws.Range("M89").Value = -6093
$ws.Range("N89").Value = -54952
$ws.Range("H105").Value = 10879.714
$ws.Range("I105").Value = 13916
$ws.Range("K105").Value = 13916
$ws.Range("M105").Value = -12169

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2894.2693
$ws.Range("I58").Value = 2482.5715
$ws.Range("K58").Value = 2482.5715
$ws.Range("M58").Value = -2279.5715
$ws.Range("H134").Value = 3342.3333
$ws.Range("I134").Value = 3013
$ws.Range("K134").Value = 9039
$ws.Range("M134").Value = -6504
$ws.Range("H136").Value = 2894.2693
$ws.Range("I136").Value = 2482.5715
$ws.Range("K136").Value = 7447.7145
$ws.Range("M136").Value = -4897.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 416.5
$ws.Range("I12").Value = 47.083332
$ws.Range("J12").Value = 859.8
$ws.Range("K12").Value = 141.249996
$ws.Range("L12").Value = 2579.4
$ws.Range("M12").Value = 31.75000399999999
$ws.Range("N12").Value = -2925.4
$ws.Range("H38").Value = 661.6087
$ws.Range("I38").Value = 226.66667
$ws.Range("K38").Value = 680.00001
$ws.Range("M38").Value = -333.00001
$ws.Range("H68").Value = 2485.1333
$ws.Range("I68").Value = 1744.75
$ws.Range("K68").Value = 5234.25
$ws.Range("M68").Value = -4423.25
$ws.Range("H71").Value = 2485.1333
$ws.Range("I71").Value = 1744.75
$ws.Range("K71").Value = 15702.75
$ws.Range("M71").Value = -11646.75
$ws.Range("H98").Value = 395.25
$ws.Range("I98").Value = 245.83333
$ws.Range("J98").Value = 843.5
$ws.Range("K98").Value = 737.49999
$ws.Range("L98").Value = 2530.5
$ws.Range("M98").Value = 760.50001
$ws.Range("N98").Value = -5526.5
$ws.Range("H107").Value = 38461950
$ws.Range("J107").Value = 83333660
$ws.Range("L107").Value = 250000980
$ws.Range("N107").Value = -250004820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 7428.5713
$ws.Range("J44").Value = 7428.5713
$ws.Range("L44").Value = 7428.5713
$ws.Range("N44").Value = -8620.5713
$ws.Range("H132").Value = 12531909
$ws.Range("I132").Value = 1745
$ws.Range("K132").Value = 5235
$ws.Range("M132").Value = -2705

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9839.352999999999
$ws.Range("I7").Value = 1478.5
$ws.Range("K7").Value = 1478.5
$ws.Range("M7").Value = -1366.5
$ws.Range("H12").Value = 250
$ws.Range("J12").Value = 250
$ws.Range("L12").Value = 250
$ws.Range("N12").Value = -590
$ws.Range("H40").Value = 4810.5557
$ws.Range("I40").Value = 3966.6667
$ws.Range("K40").Value = 3966.6667
$ws.Range("M40").Value = -3830.6667
$ws.Range("H61").Value = 10003877
$ws.Range("I61").Value = 11768385
$ws.Range("K61").Value = 11768385
$ws.Range("M61").Value = -11768183
$ws.Range("H113").Value = 10003877
$ws.Range("I113").Value = 11768385
$ws.Range("K113").Value = 11768385
$ws.Range("M113").Value = -11766215
$ws.Range("H126").Value = 9839.352999999999
$ws.Range("I126").Value = 1478.5
$ws.Range("K126").Value = 4435.5
$ws.Range("M126").Value = -1965.5
$ws.Range("H132").Value = 7851.037
$ws.Range("I132").Value = 3417.5881
$ws.Range("K132").Value = 10252.7643
$ws.Range("M132").Value = -7722.764299999999
$ws.Range("H136").Value = 4087.389
$ws.Range("I136").Value = 1843.6666
$ws.Range("K136").Value = 5530.9998
$ws.Range("M136").Value = -2980.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 15697
$ws.Range("I7").Value = 14499
$ws.Range("J7").Value = 16895
$ws.Range("K7").Value = 14499
$ws.Range("L7").Value = 16895
$ws.Range("M7").Value = -14386
$ws.Range("N7").Value = -17121
$ws.Range("H17").Value = 1801
$ws.Range("I17").Value = 1801
$ws.Range("K17").Value = 1801
$ws.Range("M17").Value = -1629
$ws.Range("H46").Value = 96197.664
$ws.Range("J46").Value = 96197.664
$ws.Range("L46").Value = 96197.664
$ws.Range("N46").Value = -96659.664
$ws.Range("H107").Value = 166667310
$ws.Range("J107").Value = 500000200
$ws.Range("L107").Value = 1500000600
$ws.Range("N107").Value = -1500004440
$ws.Range("H134").Value = 96197.664
$ws.Range("J134").Value = 96197.664
$ws.Range("L134").Value = 288592.992
$ws.Range("N134").Value = -293662.992
